$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1) - replace path-like names with slash/space -> dash/underscore
$ws.Range("B1").Value = "Government-Cadw"
$ws.Range("C1").Value = "Government-Local_Authority"
$ws.Range("D1").Value = "Government-National"
$ws.Range("E1").Value = "Government-Other"
$ws.Range("F1").Value = "Independent-English_Heritage"
$ws.Range("G1").Value = "Independent-Historic_Environment_Scotland"
$ws.Range("H1").Value = "Independent-National_Trust"
$ws.Range("I1").Value = "Independent-National_Trust_for_Scotland"
$ws.Range("J1").Value = "Independent-Not_for_profit"
$ws.Range("K1").Value = "Independent-Private"
$ws.Range("L1").Value = "Independent-Unknown"
$ws.Range("M1").Value = "University"
$ws.Range("N1").Value = "Unknown"

# Row 2 (Accredited) updated statistics
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 13.786
$ws.Range("D2").Value = 1.537
$ws.Range("E2").Value = 0.024
$ws.Range("F2").Value = 0.78
$ws.Range("G2").Value = 0.284
$ws.Range("H2").Value = 3.381
$ws.Range("I2").Value = 0.26
$ws.Range("J2").Value = 18.893
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0.024
$ws.Range("M2").Value = 1.726
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 40.695

# Row 3 (Unaccredited) updated statistics
$ws.Range("B3").Value = 0.071
$ws.Range("C3").Value = 8.016
$ws.Range("D3").Value = 0.402
$ws.Range("E3").Value = 0.213
$ws.Range("F3").Value = 0.473
$ws.Range("G3").Value = 0.213
$ws.Range("H3").Value = 0.993
$ws.Range("I3").Value = 0.378
$ws.Range("J3").Value = 22.109
$ws.Range("K3").Value = 17.758
$ws.Range("L3").Value = 5.202
$ws.Range("M3").Value = 0.875
$ws.Range("N3").Value = 2.601
$ws.Range("O3").Value = 59.304

# Row 4 (COL_TOT) updated statistics
$ws.Range("B4").Value = 0.071
$ws.Range("C4").Value = 21.802
$ws.Range("D4").Value = 1.939
$ws.Range("E4").Value = 0.237
$ws.Range("F4").Value = 1.253
$ws.Range("G4").Value = 0.497
$ws.Range("H4").Value = 4.374
$ws.Range("I4").Value = 0.638
$ws.Range("J4").Value = 41.002
$ws.Range("K4").Value = 17.758
$ws.Range("L4").Value = 5.226
$ws.Range("M4").Value = 2.601
$ws.Range("N4").Value = 2.601
$ws.Range("O4").Value = 99.999
